$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.946.41'
$ws.Range('E2').Value = '  -0.82%  '
$ws.Range('D3').Value = '3.465.94'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'591.87"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.90%  '
$ws.Range('D6').Value = "'175.27"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.73%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = "'0.583"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.90%  '
$ws.Range('D9').Value = "'0.128"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.02%  '
$ws.Range('E10').Value = '  -3.35%  '
$ws.Range('D11').Value = "'0.422"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.58%  '
$ws.Range('D12').Value = '4.060.98'
$ws.Range('E12').Value = '  -1.97%  '
$ws.Range('D13').Value = "'31.13"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +7.75%  '
$ws.Range('E14').Value = '  -0.33%  '
$ws.Range('D15').Value = '66.963.55'
$ws.Range('E15').Value = '  -0.64%  '
$ws.Range('E16').Value = '  -4.69%  '
$ws.Range('D17').Value = '3.463.43'
$ws.Range('E17').Value = '  -2.27%  '
$ws.Range('D18').Value = "'6.19"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.78%  '
$ws.Range('D19').Value = "'14.24"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.20%  '
$ws.Range('D20').Value = "'384.30"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.61%  '
$ws.Range('E21').Value = '  -2.53%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = "'1.00"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.26%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').Value = "'72.53"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.49%  '
$ws.Range('D24').Value = "'5.72"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').Value = "'0.530"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.16%  '
$ws.Range('E26').Value = '  -2.75%  '
$ws.Range('D27').Value = "'10.25"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.54%  '
$ws.Range('E28').Value = '  -2.55%  '
$ws.Range('E29').Value = '  -0.24%  '
$ws.Range('D30').Value = "'6.06"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.06%  '
$ws.Range('D31').Value = "'1.40"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.81%  '
$ws.Range('D32').Value = "'2.02"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.99%  '
$ws.Range('D33').Value = "'23.36"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.30%  '
$ws.Range('D34').Value = "'7.20"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.03%  '
$ws.Range('E35').Value = '  -2.26%  '
$ws.Range('D36').Value = "'164.34"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('D37').Value = "'0.867"
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Value = "'1.91"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.01%  '
$ws.Range('D39').Value = "'6.90"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.46%  '
$ws.Range('D40').Value = "'26.98"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.65%  '
$ws.Range('D41').Value = "'4.57"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.04%  '
$ws.Range('D42').Value = "'26.02"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.32%  '
$ws.Range('D43').Value = '2.774.60'
$ws.Range('E43').Value = '  -1.14%  '
$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D44').Value = "'0.0714"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.61%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').Value = "'2.54"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.94%  '
$ws.Range('D46').Value = "'42.12"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.05%  '
$ws.Range('E47').Value = '  -5.07%  '
$ws.Range('D48').Value = "'335.71"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.97%  '
$ws.Range('E49').Value = '  -4.38%  '
$ws.Range('D50').Value = "'32.92"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.11%  '
$ws.Range('D51').Value = "'6.32"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.59%  '
